$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shape = $s.Shapes.Item(3)
$shape.Table.ApplyStyle("{DB9610FC-84D5-4AB4-AE10-75D7F48D2EE2}")
